{"js": "// Word updated from Apache POI 4.1.0 to 5.2.3: the run properties for the\n// three header-row table cells (\"Month\", \"Savings\", \"Savings for holiday!\")\n// are rewritten as bold / non-italic / non-struck-through so they round-trip\n// through the newer POI boolean serialization (\"on\"/\"off\" instead of\n// \"true\"/\"false\"). Re-apply the (unchanged) semantic formatting explicitly\n// on each run in that header row so the run properties get (re)written.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst headerRow = table.rows.getFirst();\nconst cells = headerRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < cells.items.length; i++) {\n  const paragraphs = cells.items[i].body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < paragraphs.items.length; j++) {\n    const paragraph = paragraphs.items[j];\n    const runRange = paragraph.getRange();\n    runRange.font.set({\n      bold: true,\n      italic: false,\n      strikeThrough: false\n    });\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word updated from Apache POI 4.1.0 to 5.2.3: the run properties for the\n# three header-row table cells (\"Month\", \"Savings\", \"Savings for holiday!\")\n# are rewritten as bold / non-italic / non-struck-through so they round-trip\n# through the newer POI boolean serialization (\"on\"/\"off\" instead of\n# \"true\"/\"false\"). Re-apply the (unchanged) semantic formatting explicitly\n# on each cell's range in that header row so the run properties get\n# (re)written.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$columnCount = $table.Columns.Count\n$wdCharacter = 1\n\nfor ($col = 1; $col -le $columnCount; $col++) {\n  $cell = $table.Cell(1, $col)\n  $range = $cell.Range\n  # Trim the trailing paragraph mark so only the visible text run (not the\n  # paragraph mark) receives the font formatting.\n  $range.MoveEnd($wdCharacter, -1) | Out-Null\n  $range.Font.Bold = 1\n  $range.Font.Italic = 0\n  $range.Font.StrikeThrough = 0\n}\n"}
